$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $val) {
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

# Row 2
$ws.Cells.Item(2, 4).Value = '51.298.99'
$ws.Cells.Item(2, 5).Value = '  +2.81%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '2.746.08'
$ws.Cells.Item(3, 5).Value = '  +2.69%  '

# Row 4
$ws.Cells.Item(4, 5).Value = '  +0.08%  '

# Row 5
Set-TextCell 5 4 '115.15'
$ws.Cells.Item(5, 5).Value = '  +1.39%  '

# Row 6
Set-TextCell 6 4 '332.82'
$ws.Cells.Item(6, 5).Value = '  +2.18%  '

# Row 7
$ws.Cells.Item(7, 5).Value = '  +0.39%  '

# Row 8
Set-TextCell 8 4 '1.00'
$ws.Cells.Item(8, 5).Value = '  +0.05%  '

# Row 9
Set-TextCell 9 4 '0.573'
$ws.Cells.Item(9, 5).Value = '  +3.78%  '

# Row 10
Set-TextCell 10 4 '41.37'
$ws.Cells.Item(10, 5).Value = '  +1.14%  '

# Row 11
Set-TextCell 11 4 '20.19'
$ws.Cells.Item(11, 5).Value = '  +0.31%  '

# Row 12
Set-TextCell 12 4 '0.0827'
$ws.Cells.Item(12, 5).Value = '  +0.57%  '

# Row 13
$ws.Cells.Item(13, 5).Value = '  +2.81%  '

# Row 14
Set-TextCell 14 4 '7.65'
$ws.Cells.Item(14, 5).Value = '  +3.95%  '

# Row 15
$ws.Cells.Item(15, 4).Value = '3.177.48'
$ws.Cells.Item(15, 5).Value = '  +3.08%  '

# Row 16
$ws.Cells.Item(16, 4).Value = '2.726.13'
$ws.Cells.Item(16, 5).Value = '  +2.30%  '

# Row 17
$ws.Cells.Item(17, 5).Value = '  +1.45%  '

# Row 18
$ws.Cells.Item(18, 4).Value = '51.233.27'
$ws.Cells.Item(18, 5).Value = '  +2.83%  '

# Row 19
Set-TextCell 19 4 '13.74'
$ws.Cells.Item(19, 5).Value = '  +4.46%  '

# Row 20
Set-TextCell 20 4 '3.03'
$ws.Cells.Item(20, 5).Value = '  +4.46%  '

# Row 21
$ws.Cells.Item(21, 5).Value = '  +1.21%  '

# Row 22
$ws.Cells.Item(22, 4).Value = '0.0₃0962'
$ws.Cells.Item(22, 5).Value = '  +0.33%  '

# Row 23
Set-TextCell 23 4 '280.18'
$ws.Cells.Item(23, 5).Value = '  +0.68%  '

# Row 24
Set-TextCell 24 4 '70.13'
$ws.Cells.Item(24, 5).Value = '  -2.29%  '

# Row 25
Set-TextCell 25 4 '2.63'
$ws.Cells.Item(25, 5).Value = '  +2.29%  '

# Row 26
Set-TextCell 26 4 '26.94'
$ws.Cells.Item(26, 5).Value = '  +0.34%  '

# Row 27
$ws.Cells.Item(27, 5).Value = '  -0.01%  '

# Row 28
$ws.Cells.Item(28, 5).Value = '  +2.45%  '

# Row 29
$ws.Cells.Item(29, 5).Value = '  -0.55%  '

# Row 30
Set-TextCell 30 4 '36.03'
$ws.Cells.Item(30, 5).Value = '  -0.55%  '

# Row 31
$ws.Cells.Item(31, 5).Value = '  -0.85%  '

# Row 32
$ws.Cells.Item(32, 5).Value = '  -0.47%  '

# Row 33
$ws.Cells.Item(33, 5).Value = '  +2.47%  '

# Row 34
Set-TextCell 34 4 '0.0828'
$ws.Cells.Item(34, 5).Value = '  +2.31%  '

# Row 35
Set-TextCell 35 4 '19.64'
$ws.Cells.Item(35, 5).Value = '  +0.69%  '

# Row 36
$ws.Cells.Item(36, 5).Value = '  -0.03%  '

# Row 37
$ws.Cells.Item(37, 5).Value = '  +1.93%  '

# Row 38
Set-TextCell 38 4 '5.03'
$ws.Cells.Item(38, 5).Value = '  -0.90%  '

# Row 39
$ws.Cells.Item(39, 5).Value = '  +1.72%  '

# Row 40
$ws.Cells.Item(40, 5).Value = '  +3.58%  '

# Row 41
Set-TextCell 41 4 '23.79'
$ws.Cells.Item(41, 5).Value = '  +4.58%  '

# Row 42
$ws.Cells.Item(42, 5).Value = '  +10.50%  '

# Row 43
$ws.Cells.Item(43, 2).Value = 'WEMIXToken'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextCell 43 4 '2.30'
$ws.Cells.Item(43, 5).Value = '  +4.08%  '

# Row 44
$ws.Cells.Item(44, 2).Value = 'Stellar'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextCell 44 4 '0.113'
$ws.Cells.Item(44, 5).Value = '  +0.46%  '

# Row 45
$ws.Cells.Item(45, 2).Value = 'NEARProtocol'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextCell 45 4 '3.41'
$ws.Cells.Item(45, 5).Value = '  +3.18%  '

# Row 46
$ws.Cells.Item(46, 2).Value = 'Stacks'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextCell 46 4 '2.25'
$ws.Cells.Item(46, 5).Value = '  +10.04%  '

# Row 47
$ws.Cells.Item(47, 2).Value = 'Maker'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Cells.Item(47, 4).Value = '2.111.48'
$ws.Cells.Item(47, 5).Value = '  +0.19%  '

# Row 48
Set-TextCell 48 4 '2.29'
$ws.Cells.Item(48, 5).Value = '  +1.34%  '

# Row 49
Set-TextCell 49 4 '5.57'
$ws.Cells.Item(49, 5).Value = '  +4.00%  '

# Row 50
Set-TextCell 50 4 '9.05'
$ws.Cells.Item(50, 5).Value = '  +0.17%  '

# Row 51
$ws.Cells.Item(51, 5).Value = '  +9.88%  '

